$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "Total number of sounds in final dataset"
$ws.Range("C1").Value = 11089

# row 3
$ws.Range("C3").Value = 68
$ws.Range("D3").Value = 442
$ws.Range("F3").Value = 300

# row 5
$ws.Range("C5").Value = 45
$ws.Range("F5").Value = 122

# row 6
$ws.Range("D6").Value = 469

# row 7
$ws.Range("C7").Value = 152
$ws.Range("F7").Value = 154
$ws.Range("G7").Value = 66

# row 8
$ws.Range("C8").Value = 51
$ws.Range("D8").Value = 288
$ws.Range("F8").Value = 300

# row 9
$ws.Range("C9").Value = 52
$ws.Range("D9").Value = 242
$ws.Range("F9").Value = 294
$ws.Range("G9").Value = 23

# row 11
$ws.Range("C11").Value = 54
$ws.Range("F11").Value = 98
$ws.Range("G11").Value = 24

# row 12
$ws.Range("D12").Value = 93
$ws.Range("F12").Value = 153

# row 14
$ws.Range("D14").Value = 1307
$ws.Range("F14").Value = 300

# row 15
$ws.Range("D15").Value = 272

# row 17
$ws.Range("D17").Value = 1960
$ws.Range("F17").Value = 300

# row 19
$ws.Range("D19").Value = 1149
$ws.Range("F19").Value = 300

# row 20
$ws.Range("F20").Value = 300

# row 21
$ws.Range("D21").Value = 1659

# row 22
$ws.Range("D22").Value = 127
$ws.Range("E22").Value = 7
$ws.Range("F22").Value = 221

# row 25
$ws.Range("D25").Value = 6168
$ws.Range("E25").Value = 13
$ws.Range("F25").Value = 300

# row 26
$ws.Range("C26").Value = 49
$ws.Range("D26").Value = 2379
$ws.Range("E26").Value = 11
$ws.Range("F26").Value = 300
$ws.Range("G26").Value = 21

# row 27
$ws.Range("D27").Value = 1055
$ws.Range("F27").Value = 300

# row 29
$ws.Range("D29").Value = 128
$ws.Range("F29").Value = 180

# row 30
$ws.Range("D30").Value = 1402
$ws.Range("F30").Value = 300

# row 31
$ws.Range("C31").Value = 255
$ws.Range("D31").Value = 1516
$ws.Range("F31").Value = 300
$ws.Range("G31").Value = 110

# row 32
$ws.Range("D32").Value = 325
$ws.Range("F32").Value = 300

# row 33
$ws.Range("D33").Value = 1107
$ws.Range("F33").Value = 300

# row 35
$ws.Range("C35").Value = 65
$ws.Range("F35").Value = 251

# row 36
$ws.Range("D36").Value = 865
$ws.Range("F36").Value = 300

# row 37
$ws.Range("D37").Value = 328
$ws.Range("E37").Value = 18
$ws.Range("F37").Value = 300

# row 39
$ws.Range("D39").Value = 613
$ws.Range("F39").Value = 300

# row 40
$ws.Range("C40").Value = 77
$ws.Range("D40").Value = 132
$ws.Range("F40").Value = 209
$ws.Range("G40").Value = 33

# row 41
$ws.Range("D41").Value = 153
$ws.Range("F41").Value = 212
